$wb = $excel.ActiveWorkbook

# Remember which sheet was originally active so we can restore the
# selection at the end (inserting/renaming sheets below moves focus).
$originallyActive = $wb.Worksheets.Item(1)

# --- Step 1: Insert a new "2022-Q1" sheet right before the "总计" sheet ---
$totalSheet = $wb.Worksheets.Item(4)
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# NOTE: after Add(), the sheet reference captured in $totalSheet now tracks
# the newly inserted sheet (the collection re-targets by position), so the
# "总计" sheet must be re-fetched by name.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy the header row + A-column index cell formatting (style) from the
# "2021-Q3" sheet so the new sheet matches the existing quarterly sheets.
$srcSheet = $wb.Worksheets.Item("2021-Q3")
$srcSheet.Range("B1:H1").Copy($newSheet.Range("B1"))
$srcSheet.Range("A2:A3").Copy($newSheet.Range("A2"))

# Update the header text for the new sheet's columns.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fill in the 2022-Q1 fund holdings data (row 2). Values are entered with a
# leading apostrophe so numeric-looking strings ("007965", "0.25", ...) are
# stored as text, matching the source data (leading zeros must survive).
$newSheet.Range("B2").Value = "'007965"
$newSheet.Range("C2").Value = "民生加银品质消费股票A"
$newSheet.Range("D2").Value = "'0.25"
$newSheet.Range("E2").Value = "'88.39"
$newSheet.Range("F2").Value = "'5.13"
$newSheet.Range("G2").Value = "'0.0128"
$newSheet.Range("H2").Value = 7

# Fill in the 2022-Q1 fund holdings data (row 3).
$newSheet.Range("B3").Value = "'007966"
$newSheet.Range("C3").Value = "民生加银品质消费股票C"
$newSheet.Range("D3").Value = "'0.13"
$newSheet.Range("E3").Value = "'88.39"
$newSheet.Range("F3").Value = "'5.13"
$newSheet.Range("G3").Value = "'0.0067"
$newSheet.Range("H3").Value = 7

# --- Step 2: Insert a new top data row into the "总计" sheet for 2022-Q1 ---
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.02

# Re-apply the index-column style onto the newly inserted A2, then
# renumber the whole index column (0,1,2,3) to account for the extra row.
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# Restore the original active sheet/selection.
$originallyActive.Activate()

Write-Host "2022-Q1 sheet added; 总计 sheet updated"
